# Applies the "feat: add new changes" edit to Guia.docx
#
# Summary of the edit:
#  - A large number of paragraphs had their runs re-typed as single runs,
#    which also drops the now-stale <w:proofErr> spell-check markers that
#    wrapped individual words (Word re-creates a single run with the same
#    run properties when a paragraph's text is retyped as a whole).
#  - A brand-new step "npm i react-router-dom" was inserted at the top of
#    the "Pasos:" numbered list, pushing the rest of that list down by one
#    position (their text is unchanged other than the same run/proofErr
#    clean-up described above).
#  - A couple of sentences in the "Práctica" section were reworded
#    slightly ("... se usará como retorna el usuario ..." ->
#    "... se usará para encontrar el usuario ...").

$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

# ---------------------------------------------------------------------
# 1) Simple re-type (proofErr clean-up only, text unchanged)
# ---------------------------------------------------------------------
Set-ParaText 13 "Error boundary, componente para customizar las pantallas en caso de errores, ya sea al renderizar, en la carga de data o al ejecutar un action."
Set-ParaText 14 "Nesting routes"
Set-ParaText 21 "Hook de estado de la ruta, para saber si esta o no cargando la data."

# ---------------------------------------------------------------------
# 2) Insert the new first step of the "Pasos:" list (numId 5)
# ---------------------------------------------------------------------
$stepOne = $d.Paragraphs(23)
$stepOne.Range.InsertParagraphBefore() | Out-Null
Set-ParaText 23 "npm i react-router-dom"

# ---------------------------------------------------------------------
# 3) The rest of the "Pasos:" list shifts down by one; re-type each one
#    (text content itself is unchanged from before the insertion, only
#    proofErr/run clean-up happens).
# ---------------------------------------------------------------------
Set-ParaText 24 "Crear carpeta “pages” dentro de la carpeta “src”."
Set-ParaText 25 "Crear la carpeta “home” y el archivo index.tsx dentro de ella."
Set-ParaText 26 "Crear el componente Layout.tsx"
Set-ParaText 27 "Conectar el nuevo componente con el router para que se muestre en la ruta “/” y poner como layout el nuevo componente."
Set-ParaText 28 "Comprobar que la ruta “home” esté mostrando el componente nuevo."
Set-ParaText 29 "Agregar un botón sencillo con el texto “Ver usuarios”."
Set-ParaText 30 "Crear la carpeta “types” para definir el type de “User”."
Set-ParaText 31 "Crear la carpeta api y crear una función que simular la carga de los usuarios desde el Array, usar “defer” para estado de carga."
Set-ParaText 32 "Crear la carpeta “users” dentro de la carpeta “pages” y el archivo index.tsx con el contenido actual de App.tsx. Tomar en cuenta el uso de Suspense y Await para estado de carga."
Set-ParaText 33 "Crear la ruta “/users” como hijo de la ruta “/” y asignarle el nuevo componente Users."
Set-ParaText 34 "Definir el prop “Loader” para llamar a la función que simula la carga de usuarios."
Set-ParaText 35 "Visualizar la carga de los usuarios."
Set-ParaText 36 "Crear un nuevo componente para ErrorBoundary y definirlo como parte de la ruta “/users”."
Set-ParaText 37 "Visualizar y simular errores para comprobar los cambios."
Set-ParaText 38 "Agegar un nuevo path “/users/:state?” que tome en cuenta si solo mostrar verificados o no."

# ---------------------------------------------------------------------
# 4) "Práctica" section: proofErr clean-up plus a small rewording.
# ---------------------------------------------------------------------
Set-ParaText 40 "Consiste en realizar lo siguiente, crear una nueva ruta de la forma “users/view/{id}”, donde {id} es un parámetro obligatorio en la ruta. Este se usará para encontrar el usuario con el index en el Array. En caso de que no exista tal usuario con ese index se deberá mostrar un texto sencillo indicado que el usuario no existe. Se debe agregar un componente <a> con el texto “Ver perfil” en cada card del usuario. El index del usuario debe ser un nuevo prop de tipo numérico y opcional, debido a que solo se enviara desde la vista de lista."
Set-ParaText 41 "Tomar en cuenta el uso del prop “Loader”, la sintaxis para definir parámetros en React-router y crear una nueva función “loadUserByIndex” para devolver el usuario."
Set-ParaText 42 "Hint: Crear otro archivo view.tsx dentro de Users que use el hook userLoaderData y que use el componente UserProfile para enviar la data."

Write-Output "done"
